# ZBP_03_strategie_domacnosti.xlsx
# Adds the "16. 3. 2021" wave as a new trailing column on both sheets and
# refreshes the "aktualizace" date in the two title rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data" (percentages): new column Z, header copied from Y1 so it
# keeps the bordered/bold header style, then the 44 data rows (2-45).
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

$wsData.Range("Y1").Copy($wsData.Range("Z1"))
$wsData.Range("Z1").Value = "16. 3. 2021"

$dataZ = @{
    2  = 0.22
    3  = 0.13
    4  = 0.5
    5  = 0.3
    6  = 0.14
    7  = 0.22
    8  = 0.27
    9  = 0.18
    10 = 0.25
    11 = 0.22
    12 = 0.23
    13 = 0.37
    14 = 0.19
    15 = 0.21
    16 = 0.24
    17 = 0.19
    18 = 0.26
    19 = 0.27
    20 = 0.17
    21 = 0.16
    22 = 0.13
    23 = 0.23
    24 = 0.41
    25 = 0.45
    26 = 0.12
    27 = 0.07000000000000001
    28 = 0.13
    29 = 0.22
    30 = 0.09
    31 = 0.11
    32 = 0.13
    33 = 0.18
    34 = 0.21
    35 = 0.1
    36 = 0.14
    37 = 0.14
    38 = 0.07000000000000001
    39 = 0.25
    40 = 0.16
    41 = 0.07000000000000001
    42 = 0.05
    43 = 0.06
    44 = 0.17
    45 = 0.28
}

foreach ($row in $dataZ.Keys) {
    $wsData.Cells.Item($row, 26).Value = $dataZ[$row]
}

# Title row: bump the "aktualizace" date.
$wsData.Range("A46").Value = "Život během pandemie, Strategie domácností, % respondentů celkově a ve skupinách, aktualizace 23. 3. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR" (sample sizes): new column Y, header copied from X1,
# then the 22 data rows (2-23).
# ---------------------------------------------------------------------
$wsCount = $wb.Worksheets.Item("pocetR")

$wsCount.Range("X1").Copy($wsCount.Range("Y1"))
$wsCount.Range("Y1").Value = "16. 3. 2021"

$countY = @{
    2  = 2101
    3  = 244
    4  = 478
    5  = 1379
    6  = 1001
    7  = 187
    8  = 601
    9  = 312
    10 = 957
    11 = 174
    12 = 133
    13 = 837
    14 = 966
    15 = 721
    16 = 414
    17 = 257
    18 = 816
    19 = 615
    20 = 272
    21 = 526
    22 = 399
    23 = 235
}

foreach ($row in $countY.Keys) {
    $wsCount.Cells.Item($row, 25).Value = $countY[$row]
}

# Row 24 is the trailing "blank" row on this sheet - every other column
# (B..X) is an empty placeholder cell; mirror that for the new column Y
# by copying the empty placeholder from X24.
$wsCount.Range("X24").Copy($wsCount.Range("Y24"))

# Title row: bump the "aktualizace" date.
$wsCount.Range("A24").Value = "Život během pandemie, Strategie domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 23. 3. 2021"
